$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "MEEN30140"
$ws.Range("D37").Value = 3
$ws.Range("H37").Value = "Spr"
